$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: "home advantage" values (copy of the C5:I5 diffs) for innings 2-8
$ws.Range("A20").Value = "home advantage"
$ws.Range("C20").Value = 0.04088664000000003
$ws.Range("D20").Value = 0.05146894000000002
$ws.Range("E20").Value = 0.03722004999999995
$ws.Range("F20").Value = 0.05065118000000007
$ws.Range("G20").Value = 0.038051040000000036
$ws.Range("H20").Value = 0.037499179999999965
$ws.Range("I20").Value = 0.028128180000000003

# Row 21: average of row 20 (innings 2-8)
$ws.Range("A21").Value = "average (2-8 innings)"
$ws.Range("C21").Formula = "=AVERAGE(C20:I20)"

# Row 22: "visitor runs" values (copy of the C2:I2 values) for innings 2-8
$ws.Range("A22").Value = "visitor runs"
$ws.Range("C22").Value = 0.42001823
$ws.Range("D22").Value = 0.47471341
$ws.Range("E22").Value = 0.48881805
$ws.Range("F22").Value = 0.47699842
$ws.Range("G22").Value = 0.49805546
$ws.Range("H22").Value = 0.47474772
$ws.Range("I22").Value = 0.4656074

# Row 23: average of row 22 (innings 2-8)
$ws.Range("A23").Value = "average (2-8 innings)"
$ws.Range("C23").Formula = "=AVERAGE(C22:I22)"

# Row 24: "expected" runs in bottom 1st, computed from home advantage ratio
$ws.Range("A24").Value = """expected"" runs in bottom 1st"
$ws.Range("C24").Formula = "=B2*(1+C21/C23)"
$ws.Range("C24").Font.Bold = $true

# Update the view state to match the saved workbook
$ws.Range("A24:C24").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
